$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2026-02-10 Tuesday" "2026-02-11 Wednesday"

Replace-Text "83×27=" "29×14="
Replace-Text "47×19=" "38×44="
Replace-Text "80×19=" "21×46="
Replace-Text "31×46=" "96×85="
Replace-Text "21×83=" "75×65="
Replace-Text "49×55=" "87×51="
Replace-Text "96×45=" "68×15="
Replace-Text "27×27=" "38×80="
Replace-Text "55×97=" "43×63="
Replace-Text "47×88=" "49×19="
Replace-Text "51×13=" "51×71="
Replace-Text "40×77=" "67×63="
Replace-Text "25×14=" "81×77="
Replace-Text "17×86=" "20×52="
Replace-Text "54×43=" "38×13="
Replace-Text "77×43=" "22×89="
Replace-Text "63×50=" "94×30="
Replace-Text "46×62=" "74×62="
Replace-Text "15×78=" "27×54="
Replace-Text "63×72=" "26×70="
Replace-Text "97×60=" "36×77="
Replace-Text "79×49=" "84×32="
Replace-Text "95×14=" "87×73="
Replace-Text "35×31=" "61×36="
Replace-Text "49×66=" "29×70="
